$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.066.60"
$ws.Range("E2").Value = "  -1.65%  "
$ws.Range("D3").Value = "2.621.37"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.00"
$ws.Range("E5").Value = "  -1.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.53"
$ws.Range("E6").Value = "  -4.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.596"
$ws.Range("E7").Value = "  -2.12%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.579"
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.70"
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.31"
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0843"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.07"
$ws.Range("E13").Value = "  -3.35%  "
$ws.Range("D14").Value = "3.018.11"
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value = "2.619.26"
$ws.Range("E16").Value = "  -3.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.919"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.86"
$ws.Range("E18").Value = "  -2.80%  "
$ws.Range("D19").Value = "46.011.43"
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.77"
$ws.Range("E21").Value = "  -1.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.75"
$ws.Range("E22").Value = "  -4.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.43"
$ws.Range("E23").Value = "  +4.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "283.79"
$ws.Range("E24").Value = "  +8.63%  "
$ws.Range("E25").Value = "  -2.85%  "
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "30.20"
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.52"
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.63"
$ws.Range("E31").Value = "  -8.01%  "
$ws.Range("E32").Value = "  -4.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.23"
$ws.Range("E33").Value = "  -0.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.63"
$ws.Range("E34").Value = "  -4.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.29"
$ws.Range("E35").Value = "  -1.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "156.60"
$ws.Range("E36").Value = "  +2.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0840"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.82"
$ws.Range("E38").Value = "  -2.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.123"
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.124"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.80"
$ws.Range("E41").Value = "  -7.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.11"
$ws.Range("E42").Value = "  +1.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0327"
$ws.Range("E43").Value = "  -1.42%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.54"
$ws.Range("E44").Value = "  -5.19%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.04"
$ws.Range("E45").Value = "  -7.63%  "
$ws.Range("D46").Value = "2.109.00"
$ws.Range("E46").Value = "  +2.83%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "94.20"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "109.92"
$ws.Range("E49").Value = "  -3.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.14"
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("D51").Value = "2.870.52"
$ws.Range("E51").Value = "  -1.00%  "
